{"js": "// Insert a new \"List Bullet\" paragraph with the professor's name right\n// after the \"Docente(s) Respons\u00e1vel(eis)\" heading paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst marker = \"Docente(s) Respons\u00e1vel(eis)\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text.trim();\n  if (text === marker || text.indexOf(marker) === 0) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the 'Docente(s) Respons\u00e1vel(eis)' paragraph\");\n}\n\nconst newPara = target.insertParagraph(\"4893449 - D\u00e9bora Souza Alvim\", \"After\");\nnewPara.style = \"List Bullet\";\n\nawait context.sync();\n", "ps1": "# Insert a new \"List Bullet\" paragraph with the professor's name right\n# after the \"Docente(s) Respons\u00e1vel(eis)\" heading paragraph.\n$d = $word.ActiveDocument\n\n$marker = \"Docente(s) Respons\u00e1vel(eis)\"\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.Trim()\n    if ($t -eq $marker -or $t.StartsWith($marker)) {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the 'Docente(s) Respons\u00e1vel(eis)' paragraph\"\n}\n\n$target.Range.InsertParagraphAfter() | Out-Null\n$newPara = $target.Next()\n$newPara.Range.Text = \"4893449 - D\u00e9bora Souza Alvim\"\n$newPara.Style = \"List Bullet\"\n"}
